# Update market/profit data values on multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM)
# per scheduled market-data runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 18.333334
$ws.Range("I11").Value = 18.333334
$ws.Range("K11").Value = 18.333334
$ws.Range("M11").Value = 121.666666

# Row 40
$ws.Range("H40").Value = 1875
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1875
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1875
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2225

# Row 64
$ws.Range("H64").Value = 3200.22
$ws.Range("I64").Value = 3072.2273
$ws.Range("J64").Value = 3300.7856
$ws.Range("K64").Value = 3072.2273
$ws.Range("L64").Value = 3300.7856
$ws.Range("M64").Value = -2824.2273
$ws.Range("N64").Value = -3796.7856

# Row 67
$ws.Range("H67").Value = 3200.22
$ws.Range("I67").Value = 3072.2273
$ws.Range("J67").Value = 3300.7856
$ws.Range("K67").Value = 3072.2273
$ws.Range("L67").Value = 3300.7856
$ws.Range("M67").Value = -2214.2273
$ws.Range("N67").Value = -5016.7856

# Row 74
$ws.Range("H74").Value = 3038.8545
$ws.Range("I74").Value = 2982.8298
$ws.Range("J74").Value = 3368
$ws.Range("K74").Value = 2982.8298
$ws.Range("L74").Value = 3368
$ws.Range("M74").Value = -2046.8298
$ws.Range("N74").Value = -5240

# Row 76
$ws.Range("H76").Value = 3020.68
$ws.Range("I76").Value = 3015.85
$ws.Range("J76").Value = 3040
$ws.Range("K76").Value = 3015.85
$ws.Range("L76").Value = 3040
$ws.Range("M76").Value = -2700.85
$ws.Range("N76").Value = -3670

# Row 77
$ws.Range("H77").Value = 3038.8545
$ws.Range("I77").Value = 2982.8298
$ws.Range("J77").Value = 3368
$ws.Range("K77").Value = 14914.149
$ws.Range("L77").Value = 16840
$ws.Range("M77").Value = -10234.149
$ws.Range("N77").Value = -26200

# Row 79
$ws.Range("H79").Value = 3020.68
$ws.Range("I79").Value = 3015.85
$ws.Range("J79").Value = 3040
$ws.Range("K79").Value = 3015.85
$ws.Range("L79").Value = 3040
$ws.Range("M79").Value = -1923.85
$ws.Range("N79").Value = -5224

# Row 137
$ws.Range("H137").Value = 1814.8387
$ws.Range("I137").Value = 1577.0588
$ws.Range("J137").Value = 2103.5715
$ws.Range("K137").Value = 4731.1764
$ws.Range("L137").Value = 6310.7145
$ws.Range("M137").Value = -2181.1764
$ws.Range("N137").Value = -11410.7145

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 2482.6
$ws.Range("I63").Value = 2004.6428
$ws.Range("J63").Value = 3090.9092
$ws.Range("K63").Value = 2004.6428
$ws.Range("L63").Value = 3090.9092
$ws.Range("M63").Value = -1318.6428
$ws.Range("N63").Value = -4462.9092

# Row 66
$ws.Range("H66").Value = 2482.6
$ws.Range("I66").Value = 2004.6428
$ws.Range("J66").Value = 3090.9092
$ws.Range("K66").Value = 10023.214
$ws.Range("L66").Value = 15454.546
$ws.Range("M66").Value = -6591.214
$ws.Range("N66").Value = -22318.546

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2260.8772
$ws.Range("I105").Value = 2246.585
$ws.Range("J105").Value = 2450.25
$ws.Range("K105").Value = 2246.585
$ws.Range("L105").Value = 2450.25
$ws.Range("M105").Value = -499.585
$ws.Range("N105").Value = -5944.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2366.2104
$ws.Range("I31").Value = 1358.909
$ws.Range("J31").Value = 3751.25
$ws.Range("K31").Value = 1358.909
$ws.Range("L31").Value = 3751.25
$ws.Range("M31").Value = -1063.909
$ws.Range("N31").Value = -4341.25

# Row 34
$ws.Range("H34").Value = 2366.2104
$ws.Range("I34").Value = 1358.909
$ws.Range("J34").Value = 3751.25
$ws.Range("K34").Value = 1358.909
$ws.Range("L34").Value = 3751.25
$ws.Range("M34").Value = -1156.909
$ws.Range("N34").Value = -4155.25

# Row 62
$ws.Range("H62").Value = 2330.7144
$ws.Range("I62").Value = 2380
$ws.Range("J62").Value = 2207.5
$ws.Range("K62").Value = 2380
$ws.Range("L62").Value = 2207.5
$ws.Range("M62").Value = -1756
$ws.Range("N62").Value = -3455.5

# Row 65
$ws.Range("H65").Value = 2330.7144
$ws.Range("I65").Value = 2380
$ws.Range("J65").Value = 2207.5
$ws.Range("K65").Value = 11900
$ws.Range("L65").Value = 11037.5
$ws.Range("M65").Value = -8780
$ws.Range("N65").Value = -17277.5

# Row 141
$ws.Range("H141").Value = 35239.145
$ws.Range("J141").Value = 35239.145
$ws.Range("L141").Value = 35239.145
$ws.Range("N141").Value = -45599.145

$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 156.0625
$ws.Range("I38").Value = 72.5
$ws.Range("J38").Value = 206.2
$ws.Range("K38").Value = 217.5
$ws.Range("L38").Value = 618.5999999999999
$ws.Range("M38").Value = 129.5
$ws.Range("N38").Value = -1312.6

# Row 131
$ws.Range("H131").Value = 1194.7094
$ws.Range("J131").Value = 1095.7858
$ws.Range("L131").Value = 3287.3574
$ws.Range("N131").Value = -13367.3574

# Row 139
$ws.Range("H139").Value = 3050.4119
$ws.Range("I139").Value = 1317.4445
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 3952.3335
$ws.Range("L139").Value = 15000
$ws.Range("M139").Value = 1187.6665
$ws.Range("N139").Value = -25280

# Row 140
$ws.Range("H140").Value = 2366.4
$ws.Range("I140").Value = 778.625
$ws.Range("J140").Value = 5189.1113
$ws.Range("K140").Value = 2335.875
$ws.Range("L140").Value = 15567.3339
$ws.Range("M140").Value = 2844.125
$ws.Range("N140").Value = -25927.3339

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4771.3335
$ws.Range("I70").Value = 4742.9165
$ws.Range("J70").Value = 4794.067
$ws.Range("K70").Value = 4742.9165
$ws.Range("L70").Value = 4794.067
$ws.Range("M70").Value = -4472.9165
$ws.Range("N70").Value = -5334.067

# Row 73
$ws.Range("H73").Value = 4771.3335
$ws.Range("I73").Value = 4742.9165
$ws.Range("J73").Value = 4794.067
$ws.Range("K73").Value = 4742.9165
$ws.Range("L73").Value = 4794.067
$ws.Range("M73").Value = -3806.9165
$ws.Range("N73").Value = -6666.067

# Row 80
$ws.Range("H80").Value = 2416.3333
$ws.Range("I80").Value = 2377.7778
$ws.Range("J80").Value = 2454.889
$ws.Range("K80").Value = 2377.7778
$ws.Range("L80").Value = 2454.889
$ws.Range("M80").Value = -1379.7778
$ws.Range("N80").Value = -4450.889

# Row 83
$ws.Range("H83").Value = 2416.3333
$ws.Range("I83").Value = 2377.7778
$ws.Range("J83").Value = 2454.889
$ws.Range("K83").Value = 11888.889
$ws.Range("L83").Value = 12274.445
$ws.Range("M83").Value = -6896.888999999999
$ws.Range("N83").Value = -22258.445
